# Scheduled-runner market data refresh: updates the static
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) across the
# per-job-class leve-profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# These are plain cached numbers (not formulas), so each touched cell is
# simply overwritten with its refreshed value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 72.2
$ws.Range("I12").Value = 73.44444
$ws.Range("K12").Value = 73.44444
$ws.Range("M12").Value = 96.55556
$ws.Range("H51").Value = 19999.5
$ws.Range("J51").Value = 9999
$ws.Range("L51").Value = 9999
$ws.Range("N51").Value = -10967
$ws.Range("H86").Value = 8991.5
$ws.Range("I86").Value = 8991.5
$ws.Range("K86").Value = 8991.5
$ws.Range("M86").Value = -7868.5
$ws.Range("H88").Value = 760.125
$ws.Range("I88").Value = 635
$ws.Range("J88").Value = 835.2
$ws.Range("K88").Value = 635
$ws.Range("L88").Value = 835.2
$ws.Range("M88").Value = -229
$ws.Range("N88").Value = -1647.2
$ws.Range("H89").Value = 8991.5
$ws.Range("I89").Value = 8991.5
$ws.Range("K89").Value = 44957.5
$ws.Range("M89").Value = -39341.5
$ws.Range("H91").Value = 760.125
$ws.Range("I91").Value = 635
$ws.Range("J91").Value = 835.2
$ws.Range("K91").Value = 635
$ws.Range("L91").Value = 835.2
$ws.Range("M91").Value = 769
$ws.Range("N91").Value = -3643.2
$ws.Range("H97").Value = 3199.3572
$ws.Range("I97").Value = 814
$ws.Range("J97").Value = 3382.8462
$ws.Range("K97").Value = 2442
$ws.Range("L97").Value = 10148.5386
$ws.Range("M97").Value = -1946
$ws.Range("N97").Value = -11140.5386
$ws.Range("H98").Value = 825.64703
$ws.Range("I98").Value = 853.4545000000001
$ws.Range("J98").Value = 774.6667
$ws.Range("K98").Value = 853.4545000000001
$ws.Range("L98").Value = 774.6667
$ws.Range("M98").Value = 644.5454999999999
$ws.Range("N98").Value = -3770.6667
$ws.Range("H100").Value = 1307.6666
$ws.Range("I100").Value = 731.8
$ws.Range("K100").Value = 731.8
$ws.Range("M100").Value = -190.8
$ws.Range("H112").Value = 1377.8889
$ws.Range("I112").Value = 772.25
$ws.Range("J112").Value = 1862.4
$ws.Range("K112").Value = 2316.75
$ws.Range("L112").Value = 5587.200000000001
$ws.Range("M112").Value = -1208.75
$ws.Range("N112").Value = -7803.200000000001
$ws.Range("H122").Value = 825.64703
$ws.Range("I122").Value = 853.4545000000001
$ws.Range("J122").Value = 774.6667
$ws.Range("K122").Value = 2560.3635
$ws.Range("L122").Value = 2324.0001
$ws.Range("M122").Value = -110.3635000000004
$ws.Range("N122").Value = -7224.0001
$ws.Range("H131").Value = 6576.778
$ws.Range("I131").Value = 3148.5
$ws.Range("J131").Value = 13433.333
$ws.Range("K131").Value = 9445.5
$ws.Range("L131").Value = 40299.999
$ws.Range("M131").Value = -4405.5
$ws.Range("N131").Value = -50379.999
$ws.Range("H132").Value = 2088.3374
$ws.Range("I132").Value = 1827.6133
$ws.Range("K132").Value = 5482.8399
$ws.Range("M132").Value = -2952.8399
$ws.Range("H138").Value = 3250.775
$ws.Range("I138").Value = 2821.3333
$ws.Range("J138").Value = 3894.9375
$ws.Range("K138").Value = 8463.999899999999
$ws.Range("L138").Value = 11684.8125
$ws.Range("M138").Value = -3323.999899999999
$ws.Range("N138").Value = -21964.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1084.3513
$ws.Range("I32").Value = 674.4857
$ws.Range("K32").Value = 674.4857
$ws.Range("M32").Value = -387.4857
$ws.Range("H97").Value = 9063
$ws.Range("J97").Value = 7469.75
$ws.Range("L97").Value = 7469.75
$ws.Range("N97").Value = -8461.75
$ws.Range("H132").Value = 1524.9231
$ws.Range("I132").Value = 913.89746
$ws.Range("K132").Value = 2741.69238
$ws.Range("M132").Value = -211.69238

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3025.2104
$ws.Range("I20").Value = 1668.4286
$ws.Range("K20").Value = 1668.4286
$ws.Range("M20").Value = -1421.4286
$ws.Range("H64").Value = 2442.8948
$ws.Range("I64").Value = 4772.375
$ws.Range("J64").Value = 748.7273
$ws.Range("K64").Value = 4772.375
$ws.Range("L64").Value = 748.7273
$ws.Range("M64").Value = -4547.375
$ws.Range("N64").Value = -1198.7273
$ws.Range("H67").Value = 2442.8948
$ws.Range("I67").Value = 4772.375
$ws.Range("J67").Value = 748.7273
$ws.Range("K67").Value = 4772.375
$ws.Range("L67").Value = 748.7273
$ws.Range("M67").Value = -3992.375
$ws.Range("N67").Value = -2308.7273
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30630
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32184
$ws.Range("H99").Value = 4528.7
$ws.Range("I99").Value = 4528.7
$ws.Range("K99").Value = 4528.7
$ws.Range("M99").Value = -3030.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 125.71429
$ws.Range("I7").Value = 114.73333
$ws.Range("K7").Value = 114.73333
$ws.Range("M7").Value = -1.733329999999995
$ws.Range("H22").Value = 1473.0667
$ws.Range("J22").Value = 1485.4286
$ws.Range("L22").Value = 1485.4286
$ws.Range("N22").Value = -2185.4286
$ws.Range("H86").Value = 7930
$ws.Range("J86").Value = 8410
$ws.Range("L86").Value = 8410
$ws.Range("N86").Value = -10656
$ws.Range("H89").Value = 7930
$ws.Range("J89").Value = 8410
$ws.Range("L89").Value = 42050
$ws.Range("N89").Value = -53282
$ws.Range("H99").Value = 2058.4167
$ws.Range("I99").Value = 1779.4286
$ws.Range("K99").Value = 1779.4286
$ws.Range("M99").Value = -281.4286
$ws.Range("H126").Value = 2058.4167
$ws.Range("I126").Value = 1779.4286
$ws.Range("K126").Value = 5338.2858
$ws.Range("M126").Value = -2868.2858
$ws.Range("H132").Value = 3056.7646
$ws.Range("I132").Value = 2732.75
$ws.Range("J132").Value = 3834.4
$ws.Range("K132").Value = 8198.25
$ws.Range("L132").Value = 11503.2
$ws.Range("M132").Value = -5668.25
$ws.Range("N132").Value = -16563.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 9216.666999999999
$ws.Range("I74").Value = 5300
$ws.Range("K74").Value = 15900
$ws.Range("M74").Value = -14839
$ws.Range("H77").Value = 9216.666999999999
$ws.Range("I77").Value = 5300
$ws.Range("K77").Value = 47700
$ws.Range("M77").Value = -42396
$ws.Range("H114").Value = 3004.5
$ws.Range("I114").Value = 1013.5
$ws.Range("J114").Value = 4000
$ws.Range("K114").Value = 3040.5
$ws.Range("L114").Value = 12000
$ws.Range("M114").Value = 213.5
$ws.Range("N114").Value = -18508

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 172.8421
$ws.Range("J2").Value = 416
$ws.Range("L2").Value = 416
$ws.Range("N2").Value = -642
$ws.Range("H62").Value = 63025.668
$ws.Range("J62").Value = 65000
$ws.Range("L62").Value = 65000
$ws.Range("N62").Value = -66372
$ws.Range("H65").Value = 63025.668
$ws.Range("J65").Value = 65000
$ws.Range("L65").Value = 195000
$ws.Range("N65").Value = -201864
$ws.Range("H98").Value = 17000
$ws.Range("J98").Value = 17000
$ws.Range("L98").Value = 17000
$ws.Range("N98").Value = -22990
$ws.Range("H102").Value = 3087.2903
$ws.Range("I102").Value = 2482.6191
$ws.Range("K102").Value = 2482.6191
$ws.Range("M102").Value = -860.6190999999999
$ws.Range("H122").Value = 1835.4166
$ws.Range("I122").Value = 1463.1177
$ws.Range("K122").Value = 4389.3531
$ws.Range("M122").Value = -1939.3531

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2590.8333
$ws.Range("I7").Value = 2581.6667
$ws.Range("K7").Value = 2581.6667
$ws.Range("M7").Value = -2469.6667
$ws.Range("H16").Value = 478.66666
$ws.Range("I16").Value = 497.6
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 497.6
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = -327.6
$ws.Range("N16").Value = -440
$ws.Range("H40").Value = 1767.5333
$ws.Range("I40").Value = 1465.2858
$ws.Range("K40").Value = 1465.2858
$ws.Range("M40").Value = -1329.2858
$ws.Range("H126").Value = 2590.8333
$ws.Range("I126").Value = 2581.6667
$ws.Range("K126").Value = 7745.000100000001
$ws.Range("M126").Value = -5275.000100000001
$ws.Range("H132").Value = 2879.4827
$ws.Range("I132").Value = 1633.4445
$ws.Range("J132").Value = 3440.2
$ws.Range("K132").Value = 4900.333500000001
$ws.Range("L132").Value = 10320.6
$ws.Range("M132").Value = -2370.333500000001
$ws.Range("N132").Value = -15380.6
$ws.Range("H136").Value = 2141.8975
$ws.Range("I136").Value = 1302.55
$ws.Range("K136").Value = 3907.65
$ws.Range("M136").Value = -1357.65

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2066.8948
$ws.Range("I122").Value = 1958.0588
$ws.Range("K122").Value = 5874.1764
$ws.Range("M122").Value = -3424.1764
$ws.Range("H132").Value = 1837.1936
$ws.Range("I132").Value = 1735.3334
$ws.Range("K132").Value = 5206.0002
$ws.Range("M132").Value = -2676.0002
